# Add 2022-Q4 data
#
# 1) "总计" (summary) sheet: insert a new first data row for 2022-Q4
#    (shifting the existing quarter rows down by one row) and bump the
#    counters/row index accordingly.
# 2) Insert a brand new worksheet named "2022-Q4" right after "总计"
#    (and before "2022-Q3"), built from a duplicate of the "2022-Q3"
#    sheet (so it inherits identical styling/column layout), then
#    replace its contents with the 2022-Q4 fund holdings.
# 3) Restore "2020-Q4" (now the last tab) as the active/selected sheet,
#    since inserting/copying sheets changes the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update "总计" sheet
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

# Shift existing rows 2..9 down to 3..10, bottom-up, copying whole rows
# (A:D) so number/text formatting (style "2" on column A) travels with
# the data instead of being lost.
for ($r = 9; $r -ge 2; $r--) {
    $destRow = $r + 1
    $srcRange = $tot.Range("A" + $r + ":D" + $r)
    $dstRange = $tot.Range("A" + $destRow + ":D" + $destRow)
    $srcRange.Copy($dstRange)
}

$totalRows = @(
    @(0, "2022-Q4", 7, 0.8),
    @(1, "2022-Q3", 13, 1.1),
    @(2, "2022-Q2", 6, 0.75),
    @(3, "2022-Q1", 3, 0.12),
    @(4, "2021-Q4", 13, 1.91),
    @(5, "2021-Q3", 9, 1.12),
    @(6, "2021-Q2", 10, 2.2),
    @(7, "2021-Q1", 2, 0.86),
    @(8, "2020-Q4", 1, 0.86)
)

for ($i = 0; $i -lt $totalRows.Length; $i++) {
    $row = $i + 2
    $tot.Cells.Item($row, 1).Value = $totalRows[$i][0]
    $tot.Cells.Item($row, 2).Value = $totalRows[$i][1]
    $tot.Cells.Item($row, 3).Value = $totalRows[$i][2]
    $tot.Cells.Item($row, 4).Value = $totalRows[$i][3]
}

# ---------------------------------------------------------------------
# 2) Create the new "2022-Q4" sheet from a copy of "2022-Q3"
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q3")
$template.Copy($null, $tot)
$newSheet = $wb.Worksheets.Item("2022-Q3 (2)")
$newSheet.Name = "2022-Q4"

# The template has 13 data rows (rows 2-14); 2022-Q4 only needs 7
# (rows 2-8), so drop the extra rows.
$newSheet.Rows("9:14").Delete()

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $newSheet.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$fundRows = @(
    @(0, "000711", "嘉实医疗保健股票", "15.66", "93.63", "4.47", "0.7000", 8),
    @(1, "012704", "中银兴利稳健回报灵活配置混合A", "3.09", "81.19", "2.14", "0.0661", 9),
    @(2, "012705", "中银兴利稳健回报灵活配置混合C", "0.98", "81.19", "2.14", "0.0210", 9),
    @(3, "011886", "弘毅远方高端制造混合型发起式证券投资基金A", "0.23", "88.95", "3.61", "0.0083", 5),
    @(4, "011887", "弘毅远方高端制造混合型发起式证券投资基金C", "0.12", "88.95", "3.61", "0.0043", 5),
    @(5, "519222", "海富通欣益灵活配置混合A", "0.25", "31.65", "0.16", "0.0004", 7),
    @(6, "519221", "海富通欣益灵活配置混合C", "0.10", "31.65", "0.16", "0.0002", 7)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = $i + 2
    $rowData = $fundRows[$i]

    $newSheet.Cells.Item($r, 1).Value = $rowData[0]

    # Fund code / size / position columns are stored as text in the
    # source data (leading zeros in fund codes, fixed decimal strings
    # like "0.7000"), so force text formatting before assigning or
    # Excel auto-coerces the numeric-looking strings to numbers.
    $newSheet.Cells.Item($r, 2).NumberFormat = "@"
    $newSheet.Cells.Item($r, 2).Value = $rowData[1]

    $newSheet.Cells.Item($r, 3).Value = $rowData[2]

    $newSheet.Cells.Item($r, 4).NumberFormat = "@"
    $newSheet.Cells.Item($r, 4).Value = $rowData[3]

    $newSheet.Cells.Item($r, 5).NumberFormat = "@"
    $newSheet.Cells.Item($r, 5).Value = $rowData[4]

    $newSheet.Cells.Item($r, 6).NumberFormat = "@"
    $newSheet.Cells.Item($r, 6).Value = $rowData[5]

    $newSheet.Cells.Item($r, 7).NumberFormat = "@"
    $newSheet.Cells.Item($r, 7).Value = $rowData[6]

    $newSheet.Cells.Item($r, 8).Value = $rowData[7]
}

# ---------------------------------------------------------------------
# 3) Keep "2020-Q4" (last tab) as the selected/active sheet, matching
#    the original workbook state.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item("2020-Q4")
$lastSheet.Activate()
